# Timing issue fix - keywords, updated tc1,2 in ubc01
#
# The "CasesTab" row's query (cell B2 on sheet "startup") previously
# returned an extra `Cohort` column (coalesce(co.cohort_description, '')
# AS `Cohort`). That trailing column is removed from the Cypher query
# text so it now ends after the "Response to Treatment" column.
#
# As a side effect of the text getting one line shorter, the row height
# for row 2 shrinks from 319 to 304.5 (matching the height already used
# by rows 3 and 4), and the active selection/scroll position is reset
# to cell B2 instead of C4:E4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$b2Text = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)

MATCH (c)<--(diag:diagnosis)
 MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)
	WHERE s.clinical_study_designation IN ['UBC01'] and diag.stage_of_disease in ['T3N0M0', 'Not Applicable'] 
     OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Range("B2").Value = $b2Text

# Row 2 shrinks now that the query text has one fewer wrapped line.
$ws.Rows.Item(2).RowHeight = 304.5

# Move the active selection/view to B2 (was C4:E4).
[void]$ws.Range("B2").Select()
